# Update the "想去人数" (want-to-go count) values in column F, rows 2-5,
# on both the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2 = 59
    3 = 345
    4 = 25
    5 = 100
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
